$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IPC regiones")

# Updated data values for rows 290-297 (regional inflation data update)
$data = @{
    290 = @{ 'C'=100.80654436022658; 'D'=100.41388529125275; 'E'=100.88036396691089; 'F'=100.97719362004808; 'G'=100.63915448922428; 'H'=101.19325791651227; 'I'=100.74607051077707; 'J'=101.12377052788172; 'K'=100.87921461774297 }
    291 = @{ 'C'=100.56907758478992; 'D'=100.53358581407466; 'E'=100.47097464898768; 'F'=100.64947480416981; 'G'=100.2097946044064; 'H'=101.02950054085503; 'I'=100.36857702136638; 'J'=100.87922939445306; 'K'=100.47694717042354 }
    292 = @{ 'C'=100.69255321166003; 'D'=100.74280535688288; 'E'=100.59526536012748; 'F'=100.8244168231094; 'G'=100.35543921877249; 'H'=100.92860321930067; 'I'=100.42535155016287; 'J'=101.0518878536122; 'K'=100.64245744144543 }
    293 = @{ 'C'=100.93199241206872; 'D'=100.9538194500626; 'E'=100.94744478343137; 'F'=100.9940404647308; 'G'=100.63616576973648; 'H'=101.24822963417799; 'I'=100.57143728960307; 'J'=101.25990199261206; 'K'=100.89615451240682 }
    294 = @{ 'C'=101.22205562988458; 'D'=101.15248777192294; 'E'=101.32376152094358; 'F'=101.27728897092486; 'G'=100.77523732197081; 'H'=101.61212933586118; 'I'=100.78740155798899; 'J'=101.76885521602804; 'K'=101.21062643506568 }
    295 = @{ 'C'=101.50946460471643; 'D'=101.49321771956062; 'E'=101.70672819482401; 'F'=101.52805955517674; 'G'=101.05892013966746; 'H'=101.77680502295046; 'I'=101.13331075941143; 'J'=101.9989883102853; 'K'=101.50534173043036 }
    296 = @{ 'C'=101.90563096792266; 'D'=101.85724495061937; 'E'=102.34876321885635; 'F'=101.98347220239677; 'G'=101.52286339052453; 'H'=102.05068949102518; 'I'=101.46055220176935; 'J'=102.32261022972659; 'K'=101.95933671288608 }
    297 = @{ 'C'=101.81547901296761; 'D'=101.6651649299793; 'E'=102.53741759475174; 'F'=101.86272959674228; 'G'=101.42355415538803; 'H'=101.96226191269749; 'I'=101.27636580425535; 'J'=102.26851828736902; 'K'=101.86963146930056 }
}

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $ws.Range("$col$row").Value = $data[$row][$col]
    }
}

# Update the sheet view: scroll position and active selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 266
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I295").Select()
